# Refresh the cryptocurrency price (column D) and 1h volume-change (column
# E) figures on the active sheet, per the latest coinranking.com pull.
#
# Column D sometimes holds values that *look* numeric to Excel's automatic
# type inference (e.g. "210.53"); forcing the cell to Text first (and
# clearing the resulting direct formatting afterwards) keeps those values
# stored as text, matching the source data's original inline-string type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new Price, new Volume(1h)-without-padding)
$rows = @(
    @{ Row = 2;  D = "28.639.87";  E = "+0.96%" },
    @{ Row = 3;  D = "1.564.57";   E = "-0.44%" },
    @{ Row = 4;  D = $null;        E = "-0.10%" },
    @{ Row = 5;  D = "210.53";     E = "-0.66%" },
    @{ Row = 6;  D = $null;        E = "-0.41%" },
    @{ Row = 7;  D = $null;        E = "-0.11%" },
    @{ Row = 8;  D = "25.11";      E = "+5.98%" },
    @{ Row = 9;  D = $null;        E = "-0.16%" },
    @{ Row = 10; D = $null;        E = "-0.18%" },
    @{ Row = 11; D = $null;        E = "-0.11%" },
    @{ Row = 12; D = $null;        E = "-0.42%" },
    @{ Row = 13; D = "1.564.00";   E = "-0.48%" },
    @{ Row = 14; D = "28.654.35";  E = "+1.11%" },
    @{ Row = 15; D = "0.515";      E = "-0.23%" },
    @{ Row = 16; D = "3.65";       E = "-1.02%" },
    @{ Row = 17; D = "61.47";      E = "-0.05%" },
    @{ Row = 18; D = "231.56";     E = "+0.82%" },
    @{ Row = 19; D = $null;        E = "-0.55%" },
    @{ Row = 20; D = $null;        E = "-1.23%" },
    @{ Row = 21; D = "0.998";      E = "-0.26%" },
    @{ Row = 22; D = "3.92";       E = "-0.94%" },
    @{ Row = 23; D = "9.00";       E = "-0.30%" },
    @{ Row = 24; D = $null;        E = "+3.13%" },
    @{ Row = 25; D = "150.79";     E = "-0.23%" },
    @{ Row = 26; D = "14.81";      E = "-0.66%" },
    @{ Row = 27; D = $null;        E = "-0.13%" },
    @{ Row = 28; D = "1.00";       E = "-0.04%" },
    @{ Row = 29; D = $null;        E = "-2.09%" },
    @{ Row = 30; D = $null;        E = "-4.00%" },
    @{ Row = 31; D = $null;        E = "-1.24%" },
    @{ Row = 32; D = $null;        E = "-0.81%" },
    @{ Row = 33; D = "1.390.80";   E = "+0.34%" },
    @{ Row = 34; D = $null;        E = "-4.21%" },
    @{ Row = 35; D = $null;        E = "-3.20%" },
    @{ Row = 36; D = $null;        E = "-1.71%" },
    @{ Row = 37; D = "2.30";       E = "-2.59%" },
    @{ Row = 38; D = $null;        E = "+0.79%" },
    @{ Row = 39; D = $null;        E = "-0.73%" },
    @{ Row = 40; D = "1.94";       E = "+2.44%" },
    @{ Row = 41; D = $null;        E = "-0.19%" },
    @{ Row = 42; D = "1.00";       E = "-0.10%" },
    @{ Row = 43; D = $null;        E = "-1.17%" },
    @{ Row = 44; D = $null;        E = "-1.43%" },
    @{ Row = 45; D = "64.07";      E = "+2.91%" },
    @{ Row = 46; D = "5.27";       E = "-1.84%" },
    @{ Row = 47; D = "1.700.98";   E = "-0.51%" },
    @{ Row = 48; D = $null;        E = "-5.74%" },
    @{ Row = 49; D = "85.44";      E = "+0.32%" },
    @{ Row = 50; D = "43.22";      E = "+6.57%" },
    @{ Row = 51; D = $null;        E = "+0.49%" }
)

foreach ($entry in $rows) {
    $r = $entry.Row

    if ($null -ne $entry.D) {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $entry.D
        $dCell.ClearFormats()
    }

    $ws.Cells.Item($r, 5).Value = "  " + $entry.E + "  "
}
